$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.883.62"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.826.64"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.45"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6917"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07621"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3037"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.47"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07756"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("D12").Value = "1.824.33"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.040"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.22"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6727"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.362"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008272"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "28.870.62"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.08"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").Value = "2.075.28"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("E21").Value = "  -3.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.405"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1468"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.90"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.731"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.191"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.124"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05089"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7496"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.813"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "1.197.20"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.673"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9277"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.26"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9987"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5163"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "1.974.51"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.447"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.212"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -13.01%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.724"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "61.86"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4174"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.09%  "
